$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "63.590.89"
$ws.Range("E2").Value = "  -3.24%  "
$ws.Range("D3").Value = "2.617.18"
$ws.Range("E3").Value = "  -1.38%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").Value = "'573.60"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -4.09%  "
$ws.Range("D6").Value = "'154.65"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.38%  "
$ws.Range("E7").Value = "  +0.05%  "
$ws.Range("D8").Value = "'0.629"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.61%  "
$ws.Range("E9").Value = "  -5.65%  "
$ws.Range("D10").Value = "'5.79"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.58%  "
$ws.Range("E11").Value = "  -3.36%  "
$ws.Range("E12").Value = "  -0.37%  "
$ws.Range("D13").Value = "'28.13"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -2.03%  "
$ws.Range("D14").Value = "3.090.39"
$ws.Range("E14").Value = "  -1.27%  "
$ws.Range("E15").Value = "  -7.62%  "
$ws.Range("D16").Value = "63.518.87"
$ws.Range("E16").Value = "  -3.13%  "
$ws.Range("D17").Value = "2.622.81"
$ws.Range("E17").Value = "  -0.37%  "
$ws.Range("D18").Value = "'12.03"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -4.51%  "
$ws.Range("D19").Value = "'4.61"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -2.69%  "
$ws.Range("D20").Value = "'7.49"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.47%  "
$ws.Range("D21").Value = "'342.81"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.86%  "
$ws.Range("E22").Value = "  -0.05%  "
$ws.Range("D23").Value = "'67.07"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -2.82%  "
$ws.Range("D24").Value = "'1.73"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +3.03%  "
$ws.Range("E25").Value = "  -4.86%  "
$ws.Range("D26").Value = "'9.18"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -5.03%  "
$ws.Range("D27").Value = "'582.39"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +9.24%  "
$ws.Range("D28").Value = "'1.58"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.33%  "
$ws.Range("E29").Value = "  +0.17%  "
$ws.Range("D30").Value = "'0.160"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -2.47%  "
$ws.Range("E31").Value = "  -0.77%  "
$ws.Range("E32").Value = "  -2.72%  "
$ws.Range("D33").Value = "'1.69"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -3.97%  "
$ws.Range("D34").Value = "'6.46"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.43%  "
$ws.Range("E35").Value = "  -2.23%  "
$ws.Range("E36").Value = "  -2.56%  "
$ws.Range("D37").Value = "'19.83"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -2.70%  "
$ws.Range("E38").Value = "  -0.05%  "
$ws.Range("D39").Value = "'153.54"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -1.42%  "
$ws.Range("D40").Value = "'1.86"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -4.02%  "
$ws.Range("E41").Value = "  -0.03%  "
$ws.Range("D42").Value = "'41.22"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -2.88%  "
$ws.Range("D43").Value = "'156.60"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -2.98%  "
$ws.Range("E44").Value = "  +3.19%  "
$ws.Range("E45").Value = "  -3.18%  "
$ws.Range("D46").Value = "'0.0593"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -2.24%  "
$ws.Range("D47").Value = "'22.70"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.26%  "
$ws.Range("D48").Value = "'0.629"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.01%  "
$ws.Range("E49").Value = "  +1.91%  "
$ws.Range("E50").Value = "  -1.63%  "
$ws.Range("D51").Value = "'19.00"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -3.87%  "
